$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs | Mmp9 | Ephb2 | FAPs
$ws.Range("B2").Value = "Mmp9"
$ws.Range("C2").Value = "Ephb2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 0.001904333333333333
$ws.Range("H2").Value = 0.005713
$ws.Range("I2").Value = 0.01334809965397277
$ws.Range("J2").Value = 0.01334809965397277
$ws.Range("O2").Value = 0.921725411846598
$ws.Range("P2").Value = 0.9217254118465981
$ws.Range("Q2").Value = 0.01208538239922222
$ws.Range("R2").Value = 0.108768441593
$ws.Range("S2").Value = 0.01230328265092748
$ws.Range("T2").Value = 0.01230328265092748

# Row 3: ECs | Mmp9 | Ephb2 | MuSCs
$ws.Range("B3").Value = "Mmp9"
$ws.Range("C3").Value = "Ephb2"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("G3").Value = 0.001904333333333333
$ws.Range("H3").Value = 0.005713
$ws.Range("I3").Value = 0.01334809965397277
$ws.Range("J3").Value = 0.01334809965397277
$ws.Range("M3").Value = 0.5389353333333333
$ws.Range("N3").Value = 1.616806
$ws.Range("O3").Value = 0.07827458815340194
$ws.Range("P3").Value = 0.07827458815340194
$ws.Range("Q3").Value = 0.001026312519777778
$ws.Range("R3").Value = 0.009236812678
$ws.Range("S3").Value = 0.001044817003045286
$ws.Range("T3").Value = 0.001044817003045286

# Row 4: FAPs | Mmp9 | Ephb2 | FAPs
$ws.Range("B4").Value = "Mmp9"
$ws.Range("C4").Value = "Ephb2"
$ws.Range("D4").Value = "FAPs"
$ws.Range("I4").Value = 0.9866519003460271
$ws.Range("J4").Value = 0.9866519003460271
$ws.Range("O4").Value = 0.921725411846598
$ws.Range("P4").Value = 0.9217254118465981
$ws.Range("S4").Value = 0.9094221291956703
$ws.Range("T4").Value = 0.9094221291956704

# Row 5: FAPs | Mmp9 | Ephb2 | MuSCs
$ws.Range("B5").Value = "Mmp9"
$ws.Range("C5").Value = "Ephb2"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("I5").Value = 0.9866519003460271
$ws.Range("J5").Value = 0.9866519003460271
$ws.Range("M5").Value = 0.5389353333333333
$ws.Range("N5").Value = 1.616806
$ws.Range("O5").Value = 0.07827458815340194
$ws.Range("P5").Value = 0.07827458815340194
$ws.Range("Q5").Value = 0.07586197468088889
$ws.Range("R5").Value = 0.6827577721279999
$ws.Range("S5").Value = 0.07722977115035665
$ws.Range("T5").Value = 0.07722977115035665

# Remove the now-obsolete rows 6 and 7 (MuSCs sending cluster rows)
$ws.Rows("6:7").Delete()
